$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> Gdnf -> Ret -> ECs): updated TPM-derived values
$ws.Range("G2").Value = 0.7999296666666668
$ws.Range("H2").Value = 2.399789
$ws.Range("M2").Value = 3.21276
$ws.Range("N2").Value = 9.63828
$ws.Range("O2").Value = 0.3952107490920524
$ws.Range("P2").Value = 0.3952107490920524
$ws.Range("Q2").Value = 2.56998203588
$ws.Range("R2").Value = 23.12983832292
$ws.Range("S2").Value = 0.3952107490920524
$ws.Range("T2").Value = 0.3952107490920524

# Row 3 (MuSCs -> Gdnf -> Ret -> FAPs): updated TPM-derived values
$ws.Range("G3").Value = 0.7999296666666668
$ws.Range("H3").Value = 2.399789
$ws.Range("O3").Value = 0.4779044122124365
$ws.Range("P3").Value = 0.4779044122124365
$ws.Range("Q3").Value = 3.107723555281333
$ws.Range("R3").Value = 27.969511997532
$ws.Range("S3").Value = 0.4779044122124365
$ws.Range("T3").Value = 0.4779044122124365

# Row 4 (MuSCs -> Gdnf -> Ret -> Inflammatory-Mac becomes -> MuSCs): updated values
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.7999296666666668
$ws.Range("H4").Value = 2.399789
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.031476333333333
$ws.Range("N4").Value = 3.094429
$ws.Range("O4").Value = 0.1268848386955111
$ws.Range("P4").Value = 0.1268848386955111
$ws.Range("Q4").Value = 0.8251085194978891
$ws.Range("R4").Value = 7.425976675481001
$ws.Range("S4").Value = 0.1268848386955111
$ws.Range("T4").Value = 0.1268848386955111

# Row 5 no longer exists in the updated export; remove it entirely.
$ws.Rows(5).Delete()
